# Ajout page détail radio-tour et corrections repas
#
# Target sheet: ADMIN (room/local assignments)
# - Row 25 ("Volunteers' lounge" / "Salon des bénévoles") is corrected to
#   "Volunteers' rest room" / "Salle de repos des bénévoles" (room stays 220 A)
# - Rows 24 and 25 get a new column D formula mirroring column C
#   (same pattern already used by every other row in the table: D = C)
# - Selection moves to the newly-edited cell D25

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ADMIN")

# Correct the volunteers' room wording (room number "220 A" is unchanged)
$ws.Range("A25").Value = "Volunteers' rest room"
$ws.Range("B25").Value = "Salle de repos des bénévoles"

# Add the missing mirror formulas in column D for rows 24 and 25
$ws.Range("D24").Formula = "=C24"
$ws.Range("D25").Formula = "=C25"

# Update the sheet's selection to D25
[void]$ws.Range("D25").Select()
